# ⚡️ Datos que faltaban hasta el 10
# Rename the single sheet to "datos", add a new "metadatos" sheet after it
# describing the indicator variables, and make "metadatos" the active/
# selected sheet.

$wb = $excel.ActiveWorkbook

# --- rename the existing (and currently only) sheet ---------------------
$datos = $wb.Worksheets.Item(1)
$datos.Name = "datos"

# --- add the new "metadatos" sheet right after "datos" -------------------
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $datos)
$meta.Name = "metadatos"

# --- base font for the whole used range (non-"minor scheme" Calibri) ----
$meta.Range("A1:D7").Font.Name = "Calibri"
$meta.Range("A1:D7").Font.Size = 11
$meta.Range("A1:D7").Font.Color = 0

# --- header row -----------------------------------------------------------
$meta.Range("A1").Value = "Variables"
$meta.Range("B1").Value = "Descripción"
$meta.Range("C1").Value = "Fuente"
$meta.Range("D1").Value = "Fecha_de_extracción"

# --- anno ------------------------------------------------------------------
$meta.Range("A2").Value = "anno"
$meta.Range("B2").Value = "Año"
$meta.Range("C2").Value = "…"
$meta.Range("D2").Value = 45722

# --- codmpio ----------------------------------------------------------------
$meta.Range("A3").Value = "codmpio"
$meta.Range("B3").Value = "Código del municipio"
$meta.Range("C3").Value = "…"
$meta.Range("D3").Value = 45722

# --- numerador ---------------------------------------------------------------
$meta.Range("A4").Value = "numerador"
$meta.Range("B4").Value = "Homicidios en niños/niñas/adolescentes"
$meta.Range("C4").Value = "Instituto Nacional de Medicina Legal y Ciencias Forenses"
$meta.Range("D4").Value = 45722

# --- denominador ---------------------------------------------------------------
$meta.Range("A5").Value = "denominador"
$meta.Range("B5").Value = "Total niños/niñas/adolescentes x 100,000"
$meta.Range("C5").Value = "Departamento Administrativo Nacional de Estadística (DANE)"
$meta.Range("D5").Value = 45722

# --- homicidios (B6 stays blank) ------------------------------------------
$meta.Range("A6").Value = "homicidios"
$meta.Range("C6").Value = "Elaboración Propia"
$meta.Range("D6").Value = 45722

# row 7 is left fully blank (already the default)

# --- date formatting for the extraction-date column ------------------------
$meta.Range("D2:D6").NumberFormat = "d-mmm-yy"

# --- selection / active sheet, matching the saved UI state -----------------
$meta.Range("D2:D6").Select()
$meta.Activate()
